# Scheduled-runner update: refresh Universalis market-price snapshots
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) across the
# eight crafting-job profit sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 192.85715
$ws.Range("I5").Value = 60
$ws.Range("K5").Value = 60
$ws.Range("M5").Value = 55
$ws.Range("H6").Value = 756186.7
$ws.Range("I6").Value = 1167288.5
$ws.Range("J6").Value = 2500
$ws.Range("K6").Value = 3501865.5
$ws.Range("L6").Value = 7500
$ws.Range("M6").Value = -3501753.5
$ws.Range("N6").Value = -7724
$ws.Range("H17").Value = 1123996.2
$ws.Range("J17").Value = 1123996.2
$ws.Range("L17").Value = 3371988.6
$ws.Range("N17").Value = -3372324.6
$ws.Range("H18").Value = 1697.8
$ws.Range("I18").Value = 1647.25
$ws.Range("J18").Value = 1900
$ws.Range("K18").Value = 1647.25
$ws.Range("L18").Value = 1900
$ws.Range("M18").Value = -1363.25
$ws.Range("N18").Value = -2468
$ws.Range("H33").Value = 372
$ws.Range("J33").Value = 149
$ws.Range("L33").Value = 149
$ws.Range("N33").Value = -607
$ws.Range("H43").Value = 4667
$ws.Range("I43").Value = 4667
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4667
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4598
$ws.Range("N43").ClearContents()
$ws.Range("H116").Value = 33347446
$ws.Range("I116").Value = 62521984
$ws.Range("J116").Value = 5114.143
$ws.Range("K116").Value = 62521984
$ws.Range("L116").Value = 5114.143
$ws.Range("M116").Value = -62518542
$ws.Range("N116").Value = -11998.143
$ws.Range("H138").Value = 3461.4211
$ws.Range("I138").Value = 1725.6875
$ws.Range("J138").Value = 3812.962
$ws.Range("K138").Value = 5177.0625
$ws.Range("L138").Value = 11438.886
$ws.Range("M138").Value = -37.0625
$ws.Range("N138").Value = -21718.886

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -884
$ws.Range("H26").Value = 5296.3335
$ws.Range("I26").Value = 5296.3335
$ws.Range("K26").Value = 5296.3335
$ws.Range("M26").Value = -4966.3335
$ws.Range("H32").Value = 14679437
$ws.Range("I32").Value = 14298774
$ws.Range("J32").Value = 23815356
$ws.Range("K32").Value = 14298774
$ws.Range("L32").Value = 23815356
$ws.Range("M32").Value = -14298487
$ws.Range("N32").Value = -23815930
$ws.Range("H61").Value = 1829.123
$ws.Range("I61").Value = 1673.9286
$ws.Range("K61").Value = 1673.9286
$ws.Range("M61").Value = -1461.9286
$ws.Range("H63").Value = 3937.5
$ws.Range("I63").Value = 2375
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 2375
$ws.Range("L63").Value = 5500
$ws.Range("M63").Value = -1689
$ws.Range("N63").Value = -6872
$ws.Range("H66").Value = 3937.5
$ws.Range("I66").Value = 2375
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 11875
$ws.Range("L66").Value = 27500
$ws.Range("M66").Value = -8443
$ws.Range("N66").Value = -34364
$ws.Range("H74").Value = 2804.7666
$ws.Range("I74").Value = 2425.8
$ws.Range("K74").Value = 2425.8
$ws.Range("M74").Value = -1551.8
$ws.Range("H77").Value = 2804.7666
$ws.Range("I77").Value = 2425.8
$ws.Range("K77").Value = 12129
$ws.Range("M77").Value = -7761
$ws.Range("H136").Value = 1829.123
$ws.Range("I136").Value = 1673.9286
$ws.Range("K136").Value = 5021.7858
$ws.Range("M136").Value = -2471.7858

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 5150
$ws.Range("I11").Value = 225
$ws.Range("J11").Value = 15000
$ws.Range("K11").Value = 225
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = -85
$ws.Range("N11").Value = -15280
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 73
$ws.Range("H80").Value = 729.7059
$ws.Range("I80").Value = 606.5
$ws.Range("J80").Value = 767.61536
$ws.Range("K80").Value = 606.5
$ws.Range("L80").Value = 767.61536
$ws.Range("M80").Value = 391.5
$ws.Range("N80").Value = -2763.61536
$ws.Range("H83").Value = 729.7059
$ws.Range("I83").Value = 606.5
$ws.Range("J83").Value = 767.61536
$ws.Range("K83").Value = 3032.5
$ws.Range("L83").Value = 3838.0768
$ws.Range("M83").Value = 1959.5
$ws.Range("N83").Value = -13822.0768
$ws.Range("H94").Value = 1066.3695
$ws.Range("I94").Value = 718.375
$ws.Range("K94").Value = 718.375
$ws.Range("M94").Value = -267.375
$ws.Range("H105").Value = 2993.75
$ws.Range("I105").Value = 2678.375
$ws.Range("J105").Value = 3624.5
$ws.Range("K105").Value = 2678.375
$ws.Range("L105").Value = 3624.5
$ws.Range("M105").Value = -931.375
$ws.Range("N105").Value = -7118.5
$ws.Range("H117").Value = 69900
$ws.Range("J117").Value = 69900
$ws.Range("L117").Value = 69900
$ws.Range("N117").Value = -79078

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 3818.4
$ws.Range("I5").Value = 797
$ws.Range("J5").Value = 5832.6665
$ws.Range("K5").Value = 797
$ws.Range("L5").Value = 5832.6665
$ws.Range("M5").Value = -685
$ws.Range("N5").Value = -6056.6665
$ws.Range("H16").Value = 3166.6667
$ws.Range("I16").Value = 3000
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = -2713
$ws.Range("H31").Value = 2727.4285
$ws.Range("I31").Value = 2004.1666
$ws.Range("J31").Value = 4029.3
$ws.Range("K31").Value = 2004.1666
$ws.Range("L31").Value = 4029.3
$ws.Range("M31").Value = -1709.1666
$ws.Range("N31").Value = -4619.3
$ws.Range("H34").Value = 2727.4285
$ws.Range("I34").Value = 2004.1666
$ws.Range("J34").Value = 4029.3
$ws.Range("K34").Value = 2004.1666
$ws.Range("L34").Value = 4029.3
$ws.Range("M34").Value = -1802.1666
$ws.Range("N34").Value = -4433.3
$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 20000
$ws.Range("K60").Value = 20000
$ws.Range("M60").Value = -19489
$ws.Range("H105").Value = 1202.7142
$ws.Range("J105").Value = 310
$ws.Range("L105").Value = 310
$ws.Range("N105").Value = -3804
$ws.Range("H113").Value = 3166.6667
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 881.19354
$ws.Range("I7").Value = 764.0769
$ws.Range("J7").Value = 965.7778
$ws.Range("K7").Value = 2292.2307
$ws.Range("L7").Value = 2897.3334
$ws.Range("M7").Value = -2180.2307
$ws.Range("N7").Value = -3121.3334
$ws.Range("H124").Value = 2950
$ws.Range("I124").Value = 2950
$ws.Range("K124").Value = 8850
$ws.Range("M124").Value = -3940

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8536737
$ws.Range("J11").Value = 2805105
$ws.Range("L11").Value = 2805105
$ws.Range("N11").Value = -2805383
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744
$ws.Range("H102").Value = 1390.3914
$ws.Range("I102").Value = 1290.0476
$ws.Range("K102").Value = 1290.0476
$ws.Range("M102").Value = 331.9523999999999
$ws.Range("H120").Value = 81225.75
$ws.Range("J120").Value = 81225.75
$ws.Range("L120").Value = 81225.75
$ws.Range("N120").Value = -90901.75
$ws.Range("H122").Value = 2670.516
$ws.Range("I122").Value = 2414.4092
$ws.Range("J122").Value = 3296.5557
$ws.Range("K122").Value = 7243.2276
$ws.Range("L122").Value = 9889.667099999999
$ws.Range("M122").Value = -4793.2276
$ws.Range("N122").Value = -14789.6671
$ws.Range("H132").Value = 6832.9375
$ws.Range("I132").Value = 6936.6
$ws.Range("J132").Value = 6660.1665
$ws.Range("K132").Value = 20809.8
$ws.Range("L132").Value = 19980.4995
$ws.Range("M132").Value = -18279.8
$ws.Range("N132").Value = -25040.4995
$ws.Range("H136").Value = 60599.11
$ws.Range("I136").Value = 49000
$ws.Range("J136").Value = 62049
$ws.Range("K136").Value = 147000
$ws.Range("L136").Value = 186147
$ws.Range("M136").Value = -144450
$ws.Range("N136").Value = -191247

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9287.429
$ws.Range("I7").Value = 10670.143
$ws.Range("J7").Value = 8596.071
$ws.Range("K7").Value = 10670.143
$ws.Range("L7").Value = 8596.071
$ws.Range("M7").Value = -10558.143
$ws.Range("N7").Value = -8820.071
$ws.Range("H40").Value = 37046180
$ws.Range("I40").Value = 55562108
$ws.Range("J40").Value = 14333
$ws.Range("K40").Value = 55562108
$ws.Range("L40").Value = 14333
$ws.Range("M40").Value = -55561972
$ws.Range("N40").Value = -14605
$ws.Range("H56").Value = 8054.25
$ws.Range("I56").Value = 8054.25
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 8054.25
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -7363.25
$ws.Range("N56").ClearContents()
$ws.Range("H100").Value = 3808.4167
$ws.Range("I100").Value = 3070
$ws.Range("K100").Value = 3070
$ws.Range("M100").Value = -2529
$ws.Range("H126").Value = 9287.429
$ws.Range("I126").Value = 10670.143
$ws.Range("J126").Value = 8596.071
$ws.Range("K126").Value = 32010.429
$ws.Range("L126").Value = 25788.213
$ws.Range("M126").Value = -29540.429
$ws.Range("N126").Value = -30728.213
$ws.Range("H136").Value = 2950.1667
$ws.Range("I136").Value = 2024.75
$ws.Range("K136").Value = 6074.25
$ws.Range("M136").Value = -3524.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 20007
$ws.Range("J18").Value = 20007
$ws.Range("L18").Value = 20007
$ws.Range("N18").Value = -20353
$ws.Range("H62").Value = 4989.6
$ws.Range("I62").Value = 3599.2
$ws.Range("K62").Value = 3599.2
$ws.Range("M62").Value = -2975.2
$ws.Range("H65").Value = 4989.6
$ws.Range("I65").Value = 3599.2
$ws.Range("K65").Value = 17996
$ws.Range("M65").Value = -14876
$ws.Range("H126").Value = 3583.3438
$ws.Range("I126").Value = 3616.3215
$ws.Range("J126").Value = 3352.5
$ws.Range("K126").Value = 10848.9645
$ws.Range("L126").Value = 10057.5
$ws.Range("M126").Value = -8378.9645
$ws.Range("N126").Value = -14997.5
